$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Farsi, Pronunciation, English, HeaderStyleColumnA(1=bold+wrap, 0=plain)
$rows = @(
    @('رفتن', 'raftan', 'cleave', 1),
    @('برگشتن', 'bargashtan', 'return', 1),
    @('دادن', 'dadan', 'give', 1),
    @('برداشتن', 'bardashtan', 'take', 1),
    @('آوردن', 'avardan', 'bring', 1),
    @('دنبال __ گشتن', 'Donbale __ gashtan', 'look for', 1),
    @('پیدا کردن', 'peyda kardan', 'find', 1),
    @('گرفتن', 'gereftan', 'get ', 1),
    @('گرفتن', 'gereftan', 'receive', 1),
    @('خريدن', 'khareedan', 'buy', 1),
    @(' امتحان کردن', 'emtehan kardan', 'try', 1),
    @('شروع کردن', 'shoroo’ cardan', 'start', 1),
    @('دیگر ', 'deegar', 'stop', 1),
    @('!دیگر نخور', 'deegar nakhor!', 'stop eating', 1),
    @('تمام کردن', 'tamam kardan', 'finish', 1),
    @('ادامه دادن', 'edame dadan', 'continue', 1),
    @('بیدار شدن', 'beedar shodan', 'wake up', 1),
    @('پا شدن', 'pa shodan', 'get up', 1),
    @('خوردن', 'khordan', 'eat', 0),
    @('اتفاق افتادن', 'Et-tefagh oftadan', 'happen', 1),
    @('احساس کردن', 'ehsas kardan', 'feel', 1),
    @('درست کردن  ', 'dorost kardan', 'create/make', 1)
)

$startRow = 307
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $entry = $rows[$i]

    $farsi = $entry[0]
    $pron = $entry[1]
    $eng = $entry[2]
    $bold = $entry[3]

    $cellA = $ws.Cells.Item($r, 1)
    $cellB = $ws.Cells.Item($r, 2)
    $cellC = $ws.Cells.Item($r, 3)

    $cellA.Value = $farsi
    $cellB.Value = $pron
    $cellC.Value = $eng

    if ($bold -eq 1) {
        $cellA.WrapText = $true
        $cellA.HorizontalAlignment = 1
        $cellA.Font.Bold = $true
    }
}

$ws.Range("C329").Select()
